$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 8.307456016540527
$ws.Range("B1").Value = 6.171613216400146
$ws.Range("C1").Value = 5.016767978668213
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 3.205115556716919
